# Add "NA" values for the new column E (duplicate_image_filename) data rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

for ($row = 2; $row -le 21; $row++) {
    $ws.Range("E$row").Value = "NA"
}
